$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Donor cells with pre-existing target styles (used for format-only paste)
$donorF = $ws.Cells.Item(2, 6)   # OF_CDG date-code style (s=2)
$donorG = $ws.Cells.Item(2, 7)   # OF_DATA date style (s=3)
$donorO = $ws.Cells.Item(2, 15)  # FORNECEDOR_CDG plain style (s=1)

# ---- Row 140 ----
$donorF.Copy()
$ws.Cells.Item(140, 6).PasteSpecial(-4122)
$ws.Cells.Item(140, 6).Value = 81377
$donorG.Copy()
$ws.Cells.Item(140, 7).PasteSpecial(-4122)
$ws.Cells.Item(140, 7).Value = 46013
$ws.Cells.Item(140, 13).Value = 720
$ws.Cells.Item(140, 14).Value = 720
$ws.Cells.Item(140, 15).NumberFormat = "@"
$ws.Cells.Item(140, 15).Value = '00000000007747'
$donorO.Copy()
$ws.Cells.Item(140, 15).PasteSpecial(-4122)
$ws.Cells.Item(140, 16).Value = 'ATLANTA'

# ---- Row 141 ----
$donorF.Copy()
$ws.Cells.Item(141, 6).PasteSpecial(-4122)
$ws.Cells.Item(141, 6).Value = 81377
$donorG.Copy()
$ws.Cells.Item(141, 7).PasteSpecial(-4122)
$ws.Cells.Item(141, 7).Value = 46013
$ws.Cells.Item(141, 13).Value = 100
$ws.Cells.Item(141, 14).Value = 200
$ws.Cells.Item(141, 15).NumberFormat = "@"
$ws.Cells.Item(141, 15).Value = '00000000007747'
$donorO.Copy()
$ws.Cells.Item(141, 15).PasteSpecial(-4122)
$ws.Cells.Item(141, 16).Value = 'ATLANTA'

# ---- Row 142 ----
$ws.Cells.Item(142, 6).Value = 81378
$ws.Cells.Item(142, 7).Value = 46013
$ws.Cells.Item(142, 8).Value = 'C.05.0242'
$ws.Cells.Item(142, 9).Value = 'VENTILADOR 60 CM'
$ws.Cells.Item(142, 11).Value = 'UN'
$ws.Cells.Item(142, 13).Value = 640
$ws.Cells.Item(142, 14).Value = 640

# ---- Row 143 ----
$ws.Cells.Item(143, 8).Value = 'E.03.0150'
$ws.Cells.Item(143, 9).Value = 'BOTA  DE SEGURANÇA MSA  NOBUCK MARLUVAS CADARÇO'
$ws.Cells.Item(143, 13).Value = 190
$ws.Cells.Item(143, 14).Value = 190

# ---- Row 144 ----
$ws.Cells.Item(144, 8).Value = 'E.03.0151'
$ws.Cells.Item(144, 9).Value = 'BOTA DE SEGURANÇA  EM COURO COM SOLADO EM PU VULCAFLEX'
$ws.Cells.Item(144, 11).Value = 'PAR'
$ws.Cells.Item(144, 12).Value = 1
$ws.Cells.Item(144, 13).Value = 69.9
$ws.Cells.Item(144, 14).Value = 69.9

# ---- Row 145 ----
$ws.Cells.Item(145, 6).Value = 81378
$ws.Cells.Item(145, 7).Value = 46013
$ws.Cells.Item(145, 8).Value = 'E.04.0800'
$ws.Cells.Item(145, 9).Value = 'CORTADOR DE PISO E AZULEJO  TAM. 1,25CM'
$ws.Cells.Item(145, 13).Value = 2850
$ws.Cells.Item(145, 14).Value = 2850

# ---- Row 146 ----
$ws.Cells.Item(146, 8).Value = 'E.04.0720'
$ws.Cells.Item(146, 9).Value = 'BROXA RETANGULAR'
$ws.Cells.Item(146, 12).Value = 5
$ws.Cells.Item(146, 13).Value = 6.5
$ws.Cells.Item(146, 14).Value = 32.5

# ---- Row 147 ----
$ws.Cells.Item(147, 6).Value = 81368
$ws.Cells.Item(147, 8).Value = 'E.04.0005'
$ws.Cells.Item(147, 9).Value = 'TAMBOR DE 200L ( PARA RESERVATORIO DE AGUA.)'
$ws.Cells.Item(147, 12).Value = 1
$ws.Cells.Item(147, 13).Value = 220
$ws.Cells.Item(147, 14).Value = 220
$ws.Cells.Item(147, 15).NumberFormat = "@"
$ws.Cells.Item(147, 15).Value = '00000000008882'
$donorO.Copy()
$ws.Cells.Item(147, 15).PasteSpecial(-4122)
$ws.Cells.Item(147, 16).Value = 'GALPÃO DAS FERRAMENT'

# ---- Row 148 ----
$donorF.Copy()
$ws.Cells.Item(148, 6).PasteSpecial(-4122)
$ws.Cells.Item(148, 6).Value = 81368
$donorG.Copy()
$ws.Cells.Item(148, 7).PasteSpecial(-4122)
$ws.Cells.Item(148, 7).Value = 46010
$ws.Cells.Item(148, 8).Value = 'E.04.0646'
$ws.Cells.Item(148, 9).Value = 'LÂMINA DE SERRA P/ FERRO - STARRET'
$ws.Cells.Item(148, 12).Value = 10
$ws.Cells.Item(148, 13).Value = 13.5
$ws.Cells.Item(148, 14).Value = 135
$ws.Cells.Item(148, 15).NumberFormat = "@"
$ws.Cells.Item(148, 15).Value = '00000000008882'
$donorO.Copy()
$ws.Cells.Item(148, 15).PasteSpecial(-4122)
$ws.Cells.Item(148, 16).Value = 'GALPÃO DAS FERRAMENT'

# ---- Row 149 ----
$donorF.Copy()
$ws.Cells.Item(149, 6).PasteSpecial(-4122)
$ws.Cells.Item(149, 6).Value = 81378
$donorG.Copy()
$ws.Cells.Item(149, 7).PasteSpecial(-4122)
$ws.Cells.Item(149, 7).Value = 46013
$ws.Cells.Item(149, 8).Value = 'E.04.0066'
$ws.Cells.Item(149, 9).Value = 'MANGUEIRA FLEX PARA JARDIM DE 1/2"'
$ws.Cells.Item(149, 11).Value = 'M'
$ws.Cells.Item(149, 12).Value = 30
$ws.Cells.Item(149, 13).Value = 2.5
$ws.Cells.Item(149, 14).Value = 75
$ws.Cells.Item(149, 15).NumberFormat = "@"
$ws.Cells.Item(149, 15).Value = '00000000008882'
$donorO.Copy()
$ws.Cells.Item(149, 15).PasteSpecial(-4122)
$ws.Cells.Item(149, 16).Value = 'GALPÃO DAS FERRAMENT'

# ---- Row 150 ----
$donorF.Copy()
$ws.Cells.Item(150, 6).PasteSpecial(-4122)
$ws.Cells.Item(150, 6).Value = 81378
$donorG.Copy()
$ws.Cells.Item(150, 7).PasteSpecial(-4122)
$ws.Cells.Item(150, 7).Value = 46013
$ws.Cells.Item(150, 8).Value = 'E.04.0776'
$ws.Cells.Item(150, 9).Value = 'REGUA DE ALUMINIO 3 X 1'''' C/ 6 M'
$ws.Cells.Item(150, 11).Value = 'UN'
$ws.Cells.Item(150, 12).Value = 2
$ws.Cells.Item(150, 13).Value = 215
$ws.Cells.Item(150, 14).Value = 430
$ws.Cells.Item(150, 15).NumberFormat = "@"
$ws.Cells.Item(150, 15).Value = '00000000008882'
$donorO.Copy()
$ws.Cells.Item(150, 15).PasteSpecial(-4122)
$ws.Cells.Item(150, 16).Value = 'GALPÃO DAS FERRAMENT'

# ---- Row 151 ----
$donorF.Copy()
$ws.Cells.Item(151, 6).PasteSpecial(-4122)
$ws.Cells.Item(151, 6).Value = 81378
$donorG.Copy()
$ws.Cells.Item(151, 7).PasteSpecial(-4122)
$ws.Cells.Item(151, 7).Value = 46013
$ws.Cells.Item(151, 8).Value = 'E.04.1235'
$ws.Cells.Item(151, 9).Value = 'PISTOLA DE APLICAÇÃO (310ML )'
$ws.Cells.Item(151, 13).Value = 65
$ws.Cells.Item(151, 14).Value = 130
$ws.Cells.Item(151, 15).NumberFormat = "@"
$ws.Cells.Item(151, 15).Value = '00000000008882'
$donorO.Copy()
$ws.Cells.Item(151, 15).PasteSpecial(-4122)
$ws.Cells.Item(151, 16).Value = 'GALPÃO DAS FERRAMENT'

# ---- Row 152 ----
$donorF.Copy()
$ws.Cells.Item(152, 6).PasteSpecial(-4122)
$ws.Cells.Item(152, 6).Value = 81379
$donorG.Copy()
$ws.Cells.Item(152, 7).PasteSpecial(-4122)
$ws.Cells.Item(152, 7).Value = 46013
$ws.Cells.Item(152, 8).Value = 'K.02.0999'
$ws.Cells.Item(152, 9).Value = 'TORNEIRA DE PVC'
$ws.Cells.Item(152, 13).Value = 17
$ws.Cells.Item(152, 14).Value = 34
$ws.Cells.Item(152, 15).NumberFormat = "@"
$ws.Cells.Item(152, 15).Value = '00000000009733'
$donorO.Copy()
$ws.Cells.Item(152, 15).PasteSpecial(-4122)
$ws.Cells.Item(152, 16).Value = 'CASA PEDROSO2648864-'

# ---- Row 153 ----
$donorF.Copy()
$ws.Cells.Item(153, 6).PasteSpecial(-4122)
$ws.Cells.Item(153, 6).Value = 81379
$donorG.Copy()
$ws.Cells.Item(153, 7).PasteSpecial(-4122)
$ws.Cells.Item(153, 7).Value = 46013
$ws.Cells.Item(153, 8).Value = 'N.04.0321'
$ws.Cells.Item(153, 9).Value = 'CERÂMICA'
$ws.Cells.Item(153, 11).Value = 'M²'
$ws.Cells.Item(153, 12).Value = 1.5
$ws.Cells.Item(153, 13).Value = 30
$ws.Cells.Item(153, 14).Value = 45
$ws.Cells.Item(153, 15).NumberFormat = "@"
$ws.Cells.Item(153, 15).Value = '00000000009733'
$donorO.Copy()
$ws.Cells.Item(153, 15).PasteSpecial(-4122)
$ws.Cells.Item(153, 16).Value = 'CASA PEDROSO2648864-'

# ---- Row 154 ----
$donorF.Copy()
$ws.Cells.Item(154, 6).PasteSpecial(-4122)
$ws.Cells.Item(154, 6).Value = 81370
$donorG.Copy()
$ws.Cells.Item(154, 7).PasteSpecial(-4122)
$ws.Cells.Item(154, 7).Value = 46010
$ws.Cells.Item(154, 8).Value = 'S.10.0062'
$ws.Cells.Item(154, 9).Value = 'PU 40  FLEX ADESIVO DE POLIURETANO DE CURA RÁPIDA COR BRANCO  TUBO 310ML'
$ws.Cells.Item(154, 11).Value = 'UN'
$ws.Cells.Item(154, 12).Value = 10
$ws.Cells.Item(154, 13).Value = 10.31
$ws.Cells.Item(154, 14).Value = 103.1
$ws.Cells.Item(154, 15).NumberFormat = "@"
$ws.Cells.Item(154, 15).Value = '00000000006858'
$donorO.Copy()
$ws.Cells.Item(154, 15).PasteSpecial(-4122)
$ws.Cells.Item(154, 16).Value = 'WADY'

# ---- Row 155 ----
$donorF.Copy()
$ws.Cells.Item(155, 6).PasteSpecial(-4122)
$ws.Cells.Item(155, 6).Value = 81376
$donorG.Copy()
$ws.Cells.Item(155, 7).PasteSpecial(-4122)
$ws.Cells.Item(155, 7).Value = 46013
$ws.Cells.Item(155, 13).Value = 148.5
$ws.Cells.Item(155, 14).Value = 4455
$ws.Cells.Item(155, 15).NumberFormat = "@"
$ws.Cells.Item(155, 15).Value = '00000000008570'
$donorO.Copy()
$ws.Cells.Item(155, 15).PasteSpecial(-4122)
$ws.Cells.Item(155, 16).Value = 'ISOLIDER EPS'

Write-Host "edit complete"